$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 387, pushing the existing
# rows 387:463 down to 389:465 (dimension becomes A1:R465).
$ws.Range("387:388").Insert()

# Populate the first new row (387) with the new price-report entry.
$ws.Cells.Item(387, 1).Value = 6
$ws.Cells.Item(387, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(387, 3).Value = "Metropolitana"
$ws.Cells.Item(387, 4).Value = 45209
$ws.Cells.Item(387, 5).Value = 13
$ws.Cells.Item(387, 6).Value = 100112026
$ws.Cells.Item(387, 7).Value = "Haba"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 800
$ws.Cells.Item(387, 11).Value = 9000
$ws.Cells.Item(387, 12).Value = 10000
$ws.Cells.Item(387, 13).Value = 9475
$ws.Cells.Item(387, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(387, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(387, 16).Value = 379
$ws.Cells.Item(387, 17).Value = 25
$ws.Cells.Item(387, 18).Value = "Hortaliza"

# Populate the second new row (388) with the other new price-report entry.
$ws.Cells.Item(388, 1).Value = 6
$ws.Cells.Item(388, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(388, 3).Value = "Metropolitana"
$ws.Cells.Item(388, 4).Value = 45209
$ws.Cells.Item(388, 5).Value = 13
$ws.Cells.Item(388, 6).Value = 100112026
$ws.Cells.Item(388, 7).Value = "Haba"
$ws.Cells.Item(388, 8).Value = "Sin especificar"
$ws.Cells.Item(388, 9).Value = "Segunda"
$ws.Cells.Item(388, 10).Value = 240
$ws.Cells.Item(388, 11).Value = 6000
$ws.Cells.Item(388, 12).Value = 6000
$ws.Cells.Item(388, 13).Value = 6000
$ws.Cells.Item(388, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(388, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(388, 16).Value = 240
$ws.Cells.Item(388, 17).Value = 25
$ws.Cells.Item(388, 18).Value = "Hortaliza"
